$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the Overview paragraph's opening run right after
# "In this lab, we had" and drop a (collapsed) "_GoBack" bookmark there -
# this is what Word leaves behind at the point of the most recent edit.
# We rewrite the whole paragraph via InsertXML so we get full control over
# run/bookmark placement without picking up spurious xml:space artifacts.
# ---------------------------------------------------------------------------
$overviewPara = $d.Paragraphs.Item(6)
$overviewRange = $overviewPara.Range

$overviewXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
            </w:pPr>
            <w:r>
              <w:t>In this lab, we had</w:t>
            </w:r>
            <w:bookmarkStart w:id="101" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="101"/>
            <w:r>
              <w:t xml:space="preserve"> our first experience with conditionals and control flow in Assembly. We write a </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>while</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> loop and a short </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>if/else</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> network using equality checks and greater/less than comparisons. Also, for extra credit, we learn how to include and use another external function besides </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>printf</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> and </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>scanf</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> &#8211; </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>rand_s</w:t>
            </w:r>
            <w:r>
              <w:t>.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$overviewRange.InsertXML($overviewXml)

# ---------------------------------------------------------------------------
# Change 2: add a new bullet paragraph about the "inc" obstacle right after
# the existing "rand_s ... this caused no serious issues, though." bullet,
# and before the "Results" Heading2 paragraph.
# ---------------------------------------------------------------------------
$rngObstacle = $d.Content
$rngObstacle.Find.Execute("This caused no serious issues, though.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$obstacleEnd = $rngObstacle.End
$incTarget = $d.Range($obstacleEnd, $obstacleEnd)

$incXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Incrementing the number of guesses variable by using the </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>inc</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> command took several failed assembles, because of a forgotten pointer to the value behind the variable.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$incTarget.InsertXML($incXml)

# ---------------------------------------------------------------------------
# Change 3: the old "_GoBack" bookmark (at the very end of the document) is
# gone now that the edit point has moved to the Overview paragraph - rewrite
# the last paragraph without it.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range

$lastXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Since this was one of our first serious trips into Assembly, </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">we learned quite a bit about its syntax and semantics from this lab. Specifically, we learned about using </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>cmp</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> and various jump statements (conditional and unconditional) to implement a </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>while</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> loop and </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>if/else</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> statements.</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> As mentioned earlier, we also learned how to include and use the </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>rand_s</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> function.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$lastRange.InsertXML($lastXml)

Write-Output "edit complete"
